$d = $word.ActiveDocument

# Bump the patch version from 1.2.1 to 1.2.0 everywhere it appears
# (the document title "Version: 1.2.1 (20241204)" and the matching
# "1.2.1 (20241204)" entry heading in the Modifications changelog).
$found = $true
while ($found) {
    $found = $d.Content.Find.Execute("1.2.1", $true, $true, $false, $false, $false, `
                                      $true, 1, $false, "1.2.0", 2)
}
